$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "Uhrzeit" (time) column was removed entirely; deleting the whole
# column B shifts C:I left into B:H and keeps the existing formatting.
$ws.Columns("B").Delete()

# Fill in the previously empty Adresse (now C2) and PLZ (now F2) cells.
# Use Formula + paste-values so the text is stored as a genuine string
# (keeping the trailing space on the PLZ value) without Excel silently
# re-interpreting the PLZ text as a number or altering the cell style.
$ws.Range("C2").Formula = "=""Süsterfeldstraße 27"""
$ws.Range("C2").Copy()
$ws.Range("C2").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("F2").Formula = "=""52056 """
$ws.Range("F2").Copy()
$ws.Range("F2").PasteSpecial(-4163)  # xlPasteValues
